$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on Price cells whose new values look numeric,
# so Excel stores them as text (matching original inlineStr behavior)
# instead of silently converting to floating point numbers.
$textCells = @('D5', 'D6', 'D7', 'D8', 'D11', 'D12', 'D13', 'D14', 'D15', 'D19', 'D20', 'D21', 'D22', 'D23', 'D24', 'D25', 'D26', 'D27', 'D29', 'D30', 'D31', 'D32', 'D33', 'D34', 'D36', 'D37', 'D38', 'D39', 'D40', 'D41', 'D42', 'D43', 'D44', 'D46', 'D47', 'D48', 'D49', 'D50', 'D51')
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '89.771.14'
$ws.Range('E2').Value = '  -1.06%  '
$ws.Range('D3').Value = '3.063.88'
$ws.Range('E3').Value = '  -1.46%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = '239.69'
$ws.Range('E5').Value = '  +3.43%  '
$ws.Range('D6').Value = '614.82'
$ws.Range('E6').Value = '  -2.02%  '
$ws.Range('D7').Value = '1.12'
$ws.Range('E7').Value = '  +2.64%  '
$ws.Range('D8').Value = '0.360'
$ws.Range('E8').Value = '  -1.70%  '
$ws.Range('E9').Value = '  +0.01%  '
$ws.Range('D10').Value = '3.068.16'
$ws.Range('E10').Value = '  +6.51%  '
$ws.Range('D11').Value = '0.726'
$ws.Range('E11').Value = '  +0.44%  '
$ws.Range('D12').Value = '0.201'
$ws.Range('E12').Value = '  +2.58%  '
$ws.Range('D13').Value = '0.0000243'
$ws.Range('E13').Value = '  -1.60%  '
$ws.Range('D14').Value = '34.36'
$ws.Range('E14').Value = '  -6.08%  '
$ws.Range('D15').Value = '5.41'
$ws.Range('E15').Value = '  -1.40%  '
$ws.Range('D16').Value = '90.027.20'
$ws.Range('E16').Value = '  -0.72%  '
$ws.Range('D17').Value = '3.649.79'
$ws.Range('E17').Value = '  -1.18%  '
$ws.Range('D18').Value = '3.039.28'
$ws.Range('E18').Value = '  -2.87%  '
$ws.Range('D19').Value = '3.64'
$ws.Range('E19').Value = '  -3.99%  '
$ws.Range('D20').Value = '14.32'
$ws.Range('E20').Value = '  +1.53%  '
$ws.Range('D21').Value = '0.0000206'
$ws.Range('E21').Value = '  -0.85%  '
$ws.Range('D22').Value = '5.69'
$ws.Range('E22').Value = '  +2.65%  '
$ws.Range('D23').Value = '433.83'
$ws.Range('E23').Value = '  -1.47%  '
$ws.Range('D24').Value = '8.88'
$ws.Range('E24').Value = '  -0.13%  '
$ws.Range('B25').Value = 'NEARProtocol'
$ws.Range('C25').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D25').Value = '5.56'
$ws.Range('E25').Value = '  -2.20%  '
$ws.Range('B26').Value = 'Litecoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D26').Value = '89.89'
$ws.Range('E26').Value = '  +0.95%  '
$ws.Range('D27').Value = '11.69'
$ws.Range('E27').Value = '  -5.43%  '
$ws.Range('D28').Value = '3.243.14'
$ws.Range('E28').Value = '  -1.34%  '
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  -0.08%  '
$ws.Range('D30').Value = '0.175'
$ws.Range('E30').Value = '  +9.61%  '
$ws.Range('D31').Value = '0.235'
$ws.Range('E31').Value = '  +19.19%  '
$ws.Range('D32').Value = '9.02'
$ws.Range('E32').Value = '  -4.07%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').Value = '0.113'
$ws.Range('E33').Value = '  +30.01%  '
$ws.Range('B34').Value = 'Kaspa'
$ws.Range('C34').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D34').Value = '0.165'
$ws.Range('E34').Value = '  +9.69%  '
$ws.Range('E35').Value = '  +4.61%  '
$ws.Range('B36').Value = 'RenderToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D36').Value = '7.52'
$ws.Range('E36').Value = '  +6.87%  '
$ws.Range('B37').Value = 'MantraDAO'
$ws.Range('C37').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D37').Value = '4.20'
$ws.Range('E37').Value = '  +25.41%  '
$ws.Range('B38').Value = 'EthereumClassic'
$ws.Range('C38').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D38').Value = '25.85'
$ws.Range('E38').Value = '  -1.41%  '
$ws.Range('D39').Value = '1.88'
$ws.Range('E39').Value = '  -2.12%  '
$ws.Range('D40').Value = '477.64'
$ws.Range('E40').Value = '  -6.28%  '
$ws.Range('D41').Value = '3.50'
$ws.Range('E41').Value = '  -7.12%  '
$ws.Range('D42').Value = '1.27'
$ws.Range('E42').Value = '  -1.36%  '
$ws.Range('D43').Value = '0.411'
$ws.Range('E43').Value = '  -0.02%  '
$ws.Range('D44').Value = '22.12'
$ws.Range('E44').Value = '  -0.27%  '
$ws.Range('D46').Value = '152.87'
$ws.Range('E46').Value = '  +1.64%  '
$ws.Range('D47').Value = '1.86'
$ws.Range('E47').Value = '  -2.35%  '
$ws.Range('D48').Value = '0.675'
$ws.Range('E48').Value = '  -1.43%  '
$ws.Range('B49').Value = 'OKB'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D49').Value = '44.01'
$ws.Range('E49').Value = '  -3.01%  '
$ws.Range('B50').Value = 'ImmutableX'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D50').Value = '1.32'
$ws.Range('E50').Value = '  -1.68%  '
$ws.Range('B51').Value = 'FirstDigitalUSD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D51').Value = '0.999'
$ws.Range('E51').Value = '  +0.02%  '
